$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update weight matrix values
$ws.Range("D2").Value = 0.11
$ws.Range("D3").Value = 0.3
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 3

# Update the active selection to match the target state
$ws.Range("E10").Select()
